# Applies the weekly CompStat data refresh described in the commit message
# ("New crime data collected") to the 68th Precinct worksheet:
#   - bumps the report Volume/Number and the covered week dates
#   - refreshes the Week-to-Date / 28-Day / YTD / 2-Year crime-category grid
#     (rows 15-31), including the occasional swap between a numeric cell and a
#     literal placeholder text ("0" / "***.*") used by this report for
#     suppressed/undefined small-sample figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and reporting week -----------------------------
$ws.Range("A8").Value = "Volume 32   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# --- Fix up cell formatting where a cell switches between a numeric style --
# --- and the text-placeholder style, by copying the format+type from a    --
# --- stable donor cell elsewhere on the sheet that already has the target --
# --- representation (row 14 cells, which this edit never touches).        --
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("G14").Copy($ws.Range("D20"))
$ws.Range("H14").Copy($ws.Range("E20"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("G14").Copy($ws.Range("C28"))
$ws.Range("G14").Copy($ws.Range("D28"))
$ws.Range("H14").Copy($ws.Range("E28"))
$ws.Range("G14").Copy($ws.Range("F31"))
$ws.Range("C14").Copy($ws.Range("G31"))
$ws.Range("E14").Copy($ws.Range("H31"))

# --- Step 2: write the refreshed figures -----------------------------------
$ws.Range("M15").Value = 75
$ws.Range("D16").Value = "0"
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 15
$ws.Range("K16").Value = -44.444444444444
$ws.Range("L16").Value = -44.444444444444
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -92.753623188405
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 145.454545454545
$ws.Range("I17").Value = 56
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 55.555555555555
$ws.Range("L17").Value = 7.692307692307
$ws.Range("N17").Value = -33.333333333333
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = -16.216216216216
$ws.Range("L18").Value = -29.545454545454
$ws.Range("M18").Value = -60.256410256410
$ws.Range("N18").Value = -92.191435768262
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -7.407407407407
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 145
$ws.Range("K19").Value = -31.034482758620
$ws.Range("L19").Value = -31.034482758620
$ws.Range("M19").Value = -2.912621359223
$ws.Range("N19").Value = -37.888198757764
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 433.333333333333
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -42.372881355932
$ws.Range("L20").Value = -8.108108108108
$ws.Range("M20").Value = -41.379310344827
$ws.Range("N20").Value = -95
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 30.769230769230
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 52
$ws.Range("H21").Value = 51.923076923076
$ws.Range("I21").Value = 243
$ws.Range("J21").Value = 306
$ws.Range("K21").Value = -20.588235294117
$ws.Range("L21").Value = -22.115384615384
$ws.Range("M21").Value = -19.536423841059
$ws.Range("N21").Value = -84.220779220779
$ws.Range("F22").Value = "0"
$ws.Range("L22").Value = -33.333333333333
$ws.Range("C24").Value = 25
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -19.811320754717
$ws.Range("I24").Value = 416
$ws.Range("J24").Value = 464
$ws.Range("K24").Value = -10.344827586206
$ws.Range("L24").Value = -24.637681159420
$ws.Range("M24").Value = 3.740648379052
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -52.173913043478
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = -44.927536231884
$ws.Range("I25").Value = 204
$ws.Range("J25").Value = 292
$ws.Range("K25").Value = -30.136986301369
$ws.Range("L25").Value = -28.919860627177
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = -4.761904761904
$ws.Range("I26").Value = 124
$ws.Range("J26").Value = 102
$ws.Range("K26").Value = 21.568627450980
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 4.201680672268
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = 133.333333333333
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = "0"
$ws.Range("H31").Value = "***.*"
$ws.Range("I31").Value = 2
$ws.Range("K31").Value = -33.333333333333
$ws.Range("L31").Value = 100
